$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'25.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "'3.500"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "3LEOLEO"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "'5.091"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "'0.05631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.532"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'2.984"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "7GateTokenGT"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.8111"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "'0.8425"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9FTXTokenFTT"
$ws.Range("B11").Value = "One"
$ws.Range("C11").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D11").Value = "'0.009585"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10OneONEBestin24h"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1338"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11WazirXWRX"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.06955"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02842"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09404"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001513"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("D17").Value = "'0.006147"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.106"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17BTSETokenBTSE"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3170"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03310"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D22").Value = "'3.752"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04705"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Value = "'0.004525"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009702"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "26NitroExNTX"
$ws.Range("D28").Value = "'0.0001545"
$ws.Range("D28").Style = "Normal"
$ws.Range("D41").Value = "'0.006211"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1052"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002712"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008312"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005271"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.1800"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
